$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update M2 and N2 values, and clear O2 and P2 so they no longer hold a value
$ws.Range("M2").Value = 50
$ws.Range("N2").Value = 50
$ws.Range("O2").ClearContents()
$ws.Range("P2").ClearContents()

# Move the active selection to N5, matching the saved selection state
$ws.Range("N5").Select()
